$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing and new rows with the new IGCC netting flow data (date shifted to 2025-05-XX, new day added)
$ws.Cells.Item(2,1).Value = 45799
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 12.453
$ws.Cells.Item(3,1).Value = 45799.01041666666
$ws.Cells.Item(3,2).Value = 1.082
$ws.Cells.Item(3,3).Value = 1.91
$ws.Cells.Item(4,1).Value = 45799.02083333334
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 12.529
$ws.Cells.Item(5,1).Value = 45799.03125
$ws.Cells.Item(5,2).Value = 1.906
$ws.Cells.Item(5,3).Value = 6.075
$ws.Cells.Item(6,1).Value = 45799.04166666666
$ws.Cells.Item(6,2).Value = 3.135
$ws.Cells.Item(6,3).Value = 4.713
$ws.Cells.Item(7,1).Value = 45799.05208333334
$ws.Cells.Item(7,2).Value = 0.57
$ws.Cells.Item(7,3).Value = 2.2
$ws.Cells.Item(8,1).Value = 45799.0625
$ws.Cells.Item(8,2).Value = 0.342
$ws.Cells.Item(8,3).Value = 5.295
$ws.Cells.Item(9,1).Value = 45799.07291666666
$ws.Cells.Item(9,2).Value = 0.015
$ws.Cells.Item(9,3).Value = 7.873
$ws.Cells.Item(10,1).Value = 45799.08333333334
$ws.Cells.Item(10,2).Value = 2.208
$ws.Cells.Item(10,3).Value = 1.982
$ws.Cells.Item(11,1).Value = 45799.09375
$ws.Cells.Item(11,2).Value = 25.745
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(12,1).Value = 45799.10416666666
$ws.Cells.Item(12,2).Value = 17.831
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(13,1).Value = 45799.11458333334
$ws.Cells.Item(13,2).Value = 30.847
$ws.Cells.Item(13,3).Value = 0
$ws.Cells.Item(14,1).Value = 45799.125
$ws.Cells.Item(14,2).Value = 30.197
$ws.Cells.Item(14,3).Value = 0
$ws.Cells.Item(15,1).Value = 45799.13541666666
$ws.Cells.Item(15,2).Value = 41.48
$ws.Cells.Item(15,3).Value = 0
$ws.Cells.Item(16,1).Value = 45799.14583333334
$ws.Cells.Item(16,2).Value = 50.307
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(17,1).Value = 45799.15625
$ws.Cells.Item(17,2).Value = 55.676
$ws.Cells.Item(17,3).Value = 0
$ws.Cells.Item(18,1).Value = 45799.16666666666
$ws.Cells.Item(18,2).Value = 30.654
$ws.Cells.Item(18,3).Value = 0
$ws.Cells.Item(19,1).Value = 45799.17708333334
$ws.Cells.Item(19,2).Value = 33.482
$ws.Cells.Item(19,3).Value = 0
$ws.Cells.Item(20,1).Value = 45799.1875
$ws.Cells.Item(20,2).Value = 43.8
$ws.Cells.Item(20,3).Value = 0
$ws.Cells.Item(21,1).Value = 45799.19791666666
$ws.Cells.Item(21,2).Value = 43.485
$ws.Cells.Item(21,3).Value = 0
$ws.Cells.Item(22,1).Value = 45799.20833333334
$ws.Cells.Item(22,2).Value = 33.025
$ws.Cells.Item(22,3).Value = 0
$ws.Cells.Item(23,1).Value = 45799.21875
$ws.Cells.Item(23,2).Value = 55.785
$ws.Cells.Item(23,3).Value = 0
$ws.Cells.Item(24,1).Value = 45799.22916666666
$ws.Cells.Item(24,2).Value = 22
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(25,1).Value = 45799.23958333334
$ws.Cells.Item(25,2).Value = 21.467
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(26,1).Value = 45799.25
$ws.Cells.Item(26,2).Value = 14.287
$ws.Cells.Item(26,3).Value = 0.65
$ws.Cells.Item(27,1).Value = 45799.26041666666
$ws.Cells.Item(27,2).Value = 46.5
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(28,1).Value = 45799.27083333334
$ws.Cells.Item(28,2).Value = 15.996
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(29,1).Value = 45799.28125
$ws.Cells.Item(29,2).Value = 18.311
$ws.Cells.Item(29,3).Value = 0
$ws.Cells.Item(30,1).Value = 45799.29166666666
$ws.Cells.Item(30,2).Value = 36.791
$ws.Cells.Item(30,3).Value = 0
$ws.Cells.Item(31,1).Value = 45799.30208333334
$ws.Cells.Item(31,2).Value = 40.946
$ws.Cells.Item(31,3).Value = 0
$ws.Cells.Item(32,1).Value = 45799.3125
$ws.Cells.Item(32,2).Value = 23.144
$ws.Cells.Item(32,3).Value = 0.018
$ws.Cells.Item(33,1).Value = 45799.32291666666
$ws.Cells.Item(33,2).Value = 1.389
$ws.Cells.Item(33,3).Value = 5.063
$ws.Cells.Item(34,1).Value = 45799.33333333334
$ws.Cells.Item(34,2).Value = 28.534
$ws.Cells.Item(34,3).Value = 0.074
$ws.Cells.Item(35,1).Value = 45799.34375
$ws.Cells.Item(35,2).Value = 15.688
$ws.Cells.Item(35,3).Value = 0.343
$ws.Cells.Item(36,1).Value = 45799.35416666666
$ws.Cells.Item(36,2).Value = 0.294
$ws.Cells.Item(36,3).Value = 33.59
$ws.Cells.Item(37,1).Value = 45799.36458333334
$ws.Cells.Item(37,2).Value = 0
$ws.Cells.Item(37,3).Value = 90.681
$ws.Cells.Item(38,1).Value = 45799.375
$ws.Cells.Item(38,2).Value = 0
$ws.Cells.Item(38,3).Value = 67.001
$ws.Cells.Item(39,1).Value = 45799.38541666666
$ws.Cells.Item(39,2).Value = 10.153
$ws.Cells.Item(39,3).Value = 3.724
$ws.Cells.Item(40,1).Value = 45799.39583333334
$ws.Cells.Item(40,2).Value = 0.486
$ws.Cells.Item(40,3).Value = 14.185
$ws.Cells.Item(41,1).Value = 45799.40625
$ws.Cells.Item(41,2).Value = 0
$ws.Cells.Item(41,3).Value = 17.46
$ws.Cells.Item(42,1).Value = 45799.41666666666
$ws.Cells.Item(42,2).Value = 0
$ws.Cells.Item(42,3).Value = 57.301
$ws.Cells.Item(43,1).Value = 45799.42708333334
$ws.Cells.Item(43,2).Value = 0
$ws.Cells.Item(43,3).Value = 74.409
$ws.Cells.Item(44,1).Value = 45799.4375
$ws.Cells.Item(44,2).Value = 0
$ws.Cells.Item(44,3).Value = 28.219

# Ensure the newly added rows (42-44) use the same date/time number format as column A elsewhere
$ws.Range("A42:A44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
